# Append the next day's profit record to the tracking sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 68

# Column A holds the date as plain text (matching the existing rows), not an
# Excel date serial. Briefly marking the cell as Text ("@") before assigning
# the string keeps Excel from auto-converting "10/24/2025" into a date, and
# resetting the style back to Normal afterwards avoids leaving the cell in a
# different format than its neighbours.
$dateCell = $ws.Cells.Item($lastRow, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "10/24/2025"
$dateCell.Style = "Normal"

$ws.Cells.Item($lastRow, 2).Value = 10486.84
